# The commit swaps the two embedded DrawingML themes in this deck:
#   ppt/theme/theme1.xml  "Office Theme" (colours: Office)   -> becomes "Integral" / "Red Violet"
#   ppt/theme/theme2.xml  "Integral"     (colours: Red Violet) -> becomes "Office Theme" / "Office"
#
# theme2.xml is the theme actually wired to the slide master / presentation
# (the one PowerPoint's object model exposes live), so we drive the swap by
# rewriting its 12 theme colours from the "Red Violet" palette back to the
# stock "Office" palette via ThemeColorScheme -- the same thing a user does
# by picking a different colour variant/theme from the Design tab.

function Get-OleColor([int]$r, [int]$g, [int]$b) {
    # VBA/COM RGB() packs as 0x00BBGGRR (blue high, red low) - build that here.
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

# Index order (confirmed empirically): 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1-6, 11 hlink, 12 folHlink.
$officeColors = @(
    @(0x00, 0x00, 0x00),  # 1  dk1
    @(0xFF, 0xFF, 0xFF),  # 2  lt1
    @(0x44, 0x54, 0x6A),  # 3  dk2
    @(0xE7, 0xE6, 0xE6),  # 4  lt2
    @(0x5B, 0x9B, 0xD5),  # 5  accent1
    @(0xED, 0x7D, 0x31),  # 6  accent2
    @(0xA5, 0xA5, 0xA5),  # 7  accent3
    @(0xFF, 0xC0, 0x00),  # 8  accent4
    @(0x44, 0x72, 0xC4),  # 9  accent5
    @(0x70, 0xAD, 0x47),  # 10 accent6
    @(0x05, 0x63, 0xC1),  # 11 hlink
    @(0x95, 0x4F, 0x72)   # 12 folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $rgb = $officeColors[$i - 1]
    $tcs.Colors($i).RGB = Get-OleColor $rgb[0] $rgb[1] $rgb[2]
}
